$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A11").Value = "MB2061SS4W01-CC"
$ws.Range("B11").Value = "27-10-01-1_Pos.270"

$ws.Range("E10").Select()
